$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "datetimeFigureOut" date fields: 04/27/2012 -> 04/28/2012
#    These live on the slide master and every slide layout (not on the
#    slides themselves).
# ---------------------------------------------------------------------------
function Update-DateShapes {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "04/27/2012") {
                $shp.TextFrame.TextRange.Text = "04/28/2012"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DateShapes $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 6 ("Easter Eggs"): "50+ of one tower type" -> "20+ of one tower type"
#    Edited in two single-character steps so the run splits the same way
#    PowerPoint would while retyping the value in place.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(2)
$tr6 = $shp6.TextFrame.TextRange
$found = $tr6.Find("50+ of one tower type", 0)
$tr6.Characters($found.Start, 1).Text = "2"
$tr6b = $shp6.TextFrame.TextRange
$found2 = $tr6b.Find("0+ of one tower type", 0)
$tr6b.Characters($found2.Start, 1).Text = "0"

# ---------------------------------------------------------------------------
# 3) Slide 9 ("Music"): add a new "FL Studio" bullet (level 1) after "Audacity"
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(2)
$tr9 = $shp9.TextFrame.TextRange
$tr9.InsertAfter([char]13 + "FL Studio")
